$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.618505477905273
$ws.Range("B1").Value = 2.719385862350464
$ws.Range("C1").Value = 2.925768852233887
$ws.Range("D1").Value = 3.506209850311279
$ws.Range("E1").Value = 1.956411719322205
